# Slide 19 ("Roles del Grupo"), shape 2 ("Marcador de contenido 2"):
# the "Testin" (Testing role) paragraph description needs to be split so the
# run "es el encargado de  planificar " (with a stray double space) becomes
# two separate runs: "es el encargado de " and "planificar ".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the paragraph that contains the text to split.
$target = "es el encargado de  planificar "
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i)
    $idx = $para.Text.IndexOf($target)
    if ($idx -ge 0) {
        # Select "es el encargado de  " (19 letters/spaces + the stray extra
        # space = 20 chars) and retype it as "es el encargado de " (19
        # chars, single space). This both removes the duplicate space and
        # forces PowerPoint to split the run at that boundary, leaving the
        # remainder ("planificar ") as its own run with matching formatting.
        $selectionLength = 20
        $replacement = "es el encargado de "
        $sub = $para.Characters($idx + 1, $selectionLength)
        $sub.Text = $replacement
        break
    }
}
